# Adds a new "2022-Q4" quarter sheet to the MCD workbook and records its
# summary row on the "总计" (totals) sheet.

$wb = $excel.ActiveWorkbook

$totals = $wb.Worksheets.Item(1)     # "总计"
$template = $wb.Worksheets.Item(2)   # "2022-Q3" - used as a formatting template

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by copying "2022-Q3" (keeps styles,
#    column widths, number formats, etc.) and place it right after
#    "总计" - i.e. before the existing "2022-Q3" sheet.
# ---------------------------------------------------------------------
$template.Copy([System.Type]::Missing, $totals)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Update the fund figures for the new quarter. Columns D/E/F/G hold
# numeric-looking values stored as *text* in the source data, so a
# leading apostrophe keeps them as text instead of being coerced to
# numbers. Column H is a genuine number.
$newSheet.Range("D2").Value = "'3.86"
$newSheet.Range("E2").Value = "'94.64"
$newSheet.Range("F2").Value = "'4.43"
$newSheet.Range("G2").Value = "'0.1710"
$newSheet.Range("H2").Value = 5

$newSheet.Range("E3").Value = "'94.64"
$newSheet.Range("F3").Value = "'4.43"
$newSheet.Range("G3").Value = "'0.1267"
$newSheet.Range("H3").Value = 5

$newSheet.Range("D4").Value = "'1.00"
$newSheet.Range("E4").Value = "'94.64"
$newSheet.Range("F4").Value = "'4.43"
$newSheet.Range("G4").Value = "'0.0443"
$newSheet.Range("H4").Value = 5

# ---------------------------------------------------------------------
# 2. Insert the new quarter's summary row at the top of "总计"'s data
#    (row 2), pushing the existing rows down by one.
# ---------------------------------------------------------------------
for ($r = 9; $r -ge 2; $r--) {
    $dateVal  = $totals.Cells.Item($r, 2).Value2
    $countVal = $totals.Cells.Item($r, 3).Value2
    $valueVal = $totals.Cells.Item($r, 4).Value2
    $totals.Cells.Item($r + 1, 2).Value = $dateVal
    $totals.Cells.Item($r + 1, 3).Value = $countVal
    $totals.Cells.Item($r + 1, 4).Value = $valueVal
}

$totals.Cells.Item(2, 1).Value = 0
$totals.Cells.Item(2, 2).Value = "2022-Q4"
$totals.Cells.Item(2, 3).Value = 3
$totals.Cells.Item(2, 4).Value = 0.34

$totals.Cells.Item(10, 1).Value = 8
